# The sheet has a "query" column (A) with its output-file-name columns to the
# right (old B="dbExcel", old C="WebExcel"). Add a new "StatQuery" column
# right after A holding a second (statistics) query, shifting the old B and C
# columns one position to the right (B->C, C->D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the width of column A so the new column B can match it as closely
# as this engine allows.
$colAWidth = $ws.Columns("A").ColumnWidth

# Insert a new column at B; existing columns B and C shift right to C and D.
$ws.Columns("B").Insert()

# New header + query text for the inserted column.
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN  ['Small cell lung cancer']  OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match column A's width and the wrap-text style used by A2 for the new B2 cell.
$ws.Columns("B").ColumnWidth = $colAWidth
$ws.Range("B2").WrapText = $true

# Move the active selection to A2, matching the saved view state.
[void]$ws.Range("A2").Select()
